$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column C for all data rows (2-62) from 45207 to 45208
$ws.Range("C2:C62").Value = 45208

# Update the hyperlink formulas in row 2 to point to Logging_0484 instead of Logging_ESKILSTUNA
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/artfynd/A 31987-2023.xlsx", "A 31987-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/kartor/A 31987-2023.png", "A 31987-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/klagomål/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/klagomålsmail/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/tillsyn/A 31987-2023.docx", "A 31987-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0484/tillsynsmail/A 31987-2023.docx", "A 31987-2023")'
